$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '30.199.41'
$ws.Range("D3").Value = '2.080.17'
$ws.Range("E3").Value = '  -1.73%  '
$ws.Range("E4").Value = '  -0.10%  '
$ws.Range("D5").Value = '338.83'
$ws.Range("E5").Value = '  -2.48%  '
$ws.Range("D6").Value = '1.006'
$ws.Range("E6").Value = '  -0.07%  '
$ws.Range("D7").Value = '0.5271'
$ws.Range("E7").Value = '  +1.49%  '
$ws.Range("D8").Value = '0.4363'
$ws.Range("E8").Value = '  -1.94%  '
$ws.Range("D9").Value = '54.85'
$ws.Range("E9").Value = '  +1.36%  '
$ws.Range("D10").Value = '0.09346'
$ws.Range("E10").Value = '  +0.09%  '
$ws.Range("D11").Value = '1.171'
$ws.Range("E11").Value = '  -0.77%  '
$ws.Range("D12").Value = '24.45'
$ws.Range("E12").Value = '  -2.71%  '
$ws.Range("B13").Value = 'Chainlink'
$ws.Range("C13").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D13").Value = '8.465'
$ws.Range("E13").Value = '  +0.73%  '
$ws.Range("B14").Value = 'WrappedEther'
$ws.Range("C14").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D14").Value = '2.090.65'
$ws.Range("E14").Value = '  +0.48%  '
$ws.Range("D15").Value = '6.845'
$ws.Range("E15").Value = '  -0.08%  '
$ws.Range("D16").Value = '101.02'
$ws.Range("E16").Value = '  -1.23%  '
$ws.Range("D17").Value = '0.00001158'
$ws.Range("E17").Value = '  -0.42%  '
$ws.Range("E18").Value = '  -0.06%  '
$ws.Range("D19").Value = '20.91'
$ws.Range("E19").Value = '  -2.64%  '
$ws.Range("D20").Value = '0.06714'
$ws.Range("E20").Value = '  +0.63%  '
$ws.Range("D21").Value = '6.309'
$ws.Range("E21").Value = '  +0.21%  '
$ws.Range("D22").Value = '1.004'
$ws.Range("E22").Value = '  -0.22%  '
$ws.Range("D23").Value = '30.209.54'
$ws.Range("E23").Value = '  +1.01%  '
$ws.Range("D24").Value = '12.39'
$ws.Range("E24").Value = '  -2.46%  '
$ws.Range("D25").Value = '2.316'
$ws.Range("E25").Value = '  -0.53%  '
$ws.Range("D26").Value = '21.75'
$ws.Range("E26").Value = '  -1.59%  '
$ws.Range("B27").Value = 'InternetComputer(DFINITY)'
$ws.Range("C27").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D27").Value = '6.834'
$ws.Range("E27").Value = '  +5.93%  '
$ws.Range("B28").Value = 'Monero'
$ws.Range("C28").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D28").Value = '162.30'
$ws.Range("E28").Value = '  +0.08%  '
$ws.Range("D29").Value = '2.482'
$ws.Range("E29").Value = '  -3.33%  '
$ws.Range("D30").Value = '133.39'
$ws.Range("E30").Value = '  -0.38%  '
$ws.Range("D31").Value = '1.125'
$ws.Range("E31").Value = '  -2.34%  '
$ws.Range("B32").Value = 'ARBITRUM'
$ws.Range("C32").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D32").Value = '1.660'
$ws.Range("E32").Value = '  -7.29%  '
$ws.Range("B33").Value = 'Stellar'
$ws.Range("C33").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D33").Value = '0.1047'
$ws.Range("E33").Value = '  -0.82%  '
$ws.Range("D34").Value = '6.243'
$ws.Range("E34").Value = '  -0.01%  '
$ws.Range("D35").Value = '3.912'
$ws.Range("E35").Value = '  -1.52%  '
$ws.Range("D36").Value = '0.02601'
$ws.Range("E36").Value = '  +0.19%  '
$ws.Range("D37").Value = '9.867'
$ws.Range("E37").Value = '  -9.00%  '
$ws.Range("D38").Value = '0.06725'
$ws.Range("E38").Value = '  -1.29%  '
$ws.Range("B39").Value = 'Aptos'
$ws.Range("C39").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D39").Value = '12.54'
$ws.Range("E39").Value = '  -1.10%  '
$ws.Range("D40").Value = '0.6943'
$ws.Range("E40").Value = '  -1.32%  '
$ws.Range("B41").Value = 'TrustWalletToken'
$ws.Range("C41").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D41").Value = '1.342'
$ws.Range("E41").Value = '  +0.80%  '
$ws.Range("E42").Value = '  -2.10%  '
$ws.Range("D43").Value = '0.6726'
$ws.Range("E43").Value = '  -1.37%  '
$ws.Range("E44").Value = '  +0.88%  '
$ws.Range("D45").Value = '14.24'
$ws.Range("E45").Value = '  -1.92%  '
$ws.Range("E46").Value = '  +0.06%  '
$ws.Range("D47").Value = '1.296'
$ws.Range("E47").Value = '  +5.05%  '
$ws.Range("D48").Value = '3.624'
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("B49").Value = 'ThetaToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D49").Value = '1.209'
$ws.Range("E49").Value = '  +2.23%  '
$ws.Range("B50").Value = 'BabyDogeCoin'
$ws.Range("C50").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D50").Value = '0.00000000341'
$ws.Range("E50").Value = '  -4.63%  '
$ws.Range("D51").Value = '1.210'
$ws.Range("E51").Value = '  -0.91%  '
